{"js": "// Applies the \"Bug Fixes as per trace files\" text corrections to the shell-lab\n// assignment document.  All edits are in-place text fixes (added/removed\n// spaces, stray characters, a leftover literal \"&nbsp;\" artifact) plus a\n// run merge.  We locate each change with a unique, literal search string and\n// replace just that substring, so paragraph/run formatting (rFonts/sz) is\n// left untouched everywhere except the one spot where the source OOXML also\n// merges two runs into one (handled explicitly below).\n\nconst body = context.document.body;\n\n// Simple helper: find a unique literal substring in the document and\n// replace it with new text, preserving the surrounding run formatting.\nasync function replaceOnce(searchText, replaceText) {\n  const results = body.search(searchText, { matchCase: true, matchWildcards: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"expected exactly 1 match for \" + JSON.stringify(searchText) +\n      \" but found \" + results.items.length\n    );\n  }\n\n  results.items[0].insertText(replaceText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) \"...this initial child process).tsh need not support...\"\n//    -> \"...this initial child process).t sh need not support...\"\nawait replaceOnce(\n  \"child process).tsh need not\",\n  \"child process).t sh need not\"\n);\n\n// 2) \"...run the job in the background.Otherwise, it should run the job in\n//    the foreground.Each job can be identified...\" -> add a space after each\n//    of those two sentence-ending periods.\nawait replaceOnce(\n  \"in the background.Otherwise, it should run the job in the foreground.Each job can\",\n  \"in the background. Otherwise, it should run the job in the foreground. Each job can\"\n);\n\n// 3) \"...manipulating the joblist.)tsh should support the following...\"\n//    -> \"...manipulating the joblist. tsh should support the following...\"\n//    (drops the stray close-paren and adds a space)\nawait replaceOnce(\n  \"manipulating the joblist.)tsh should support\",\n  \"manipulating the joblist. tsh should support\"\n);\n\n// 4) \"...argument can be either a PID or a JID.The fg <job> command...\"\n//    -> \"...argument can be either a PID or a JID. The fg <job> command...\"\nawait replaceOnce(\n  \"argument can be either a PID or a JID.The fg\",\n  \"argument can be either a PID or a JID. The fg\"\n);\n\n// 5) \"...check your work.Reference solution. The Linux executable...\"\n//    -> \"...check your work. Reference solution. The Linux executable...\"\nawait replaceOnce(\n  \"check your work.Reference solution.\",\n  \"check your work. Reference solution.\"\n);\n\n// 6) \"Usage: sdriver.pl [-hv] -t <trace> -s <shellprog> -a <args>&nbsp;\"\n//    -> drop the trailing stray literal \"&nbsp;\" artifact.\nawait replaceOnce(\n  \"-a <args>&nbsp;\",\n  \"-a <args>\"\n);\n\n// 7) The \"Options:\" paragraph that starts with a run containing \"  \" and a\n//    second run containing the stray literal \"&nbsp;   -v  &nbsp;         \"\n//    prefix before \"Be more verbose ...\". The source edit deletes the first\n//    run entirely and strips the \"&nbsp;   -v  &nbsp;         \" prefix from\n//    the second run, leaving a single run. Replace the *entire paragraph's*\n//    text (spanning both runs) with the cleaned-up text in one shot, using\n//    the paragraph's own range (rather than the search-hit range) so the\n//    paragraph collapses to exactly one run without picking up a stray\n//    xml:space=\"preserve\" marker, matching the target OOXML exactly.\n{\n  const marker = body.search(\"&nbsp;   -v  &nbsp;         Be more verbose\", {\n    matchCase: true,\n    matchWildcards: false\n  });\n  marker.load(\"items\");\n  await context.sync();\n\n  if (marker.items.length !== 1) {\n    throw new Error(\n      \"expected exactly 1 match for the verbose-options paragraph but found \" +\n      marker.items.length\n    );\n  }\n\n  const optionsParagraph = marker.items[0].paragraphs.getFirst();\n  const fullRange = optionsParagraph.getRange();\n  fullRange.insertText(\n    \"Be more verbose  -t <trace>    Trace file  -s <shell>    Shell program to test -a <args> Shell arguments      -g            Generate output for autograder\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n", "ps1": "# Applies the \"Bug Fixes as per trace files\" text corrections to the\n# shell-lab assignment document. All edits are in-place text fixes (added\n# or removed spaces, a stray extra character, a leftover literal \"&nbsp;\"\n# artifact) plus a two-run merge in the \"Options:\" paragraph. Each fix is\n# applied with Find/Replace against a short, unique literal substring so\n# that surrounding run formatting (rFonts/sz) is left untouched.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Once($FindText, $ReplaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $FindText\n    $find.Replacement.Text = $ReplaceText\n\n    # wdFindWrap=1, wdReplaceAll=2 -> replace every (here: the single unique) match.\n    $ok = $find.Execute($FindText, $false, $false, $false, $false, $false, $true, 1, $false, $ReplaceText, 2)\n    if (-not $ok) {\n        throw \"Find/Replace failed for: $FindText\"\n    }\n}\n\n# 1) \"...this initial child process).tsh need not support...\"\n#    -> \"...this initial child process).t sh need not support...\"\nReplace-Once \"child process).tsh need not\" \"child process).t sh need not\"\n\n# 2) \"...run the job in the background.Otherwise, it should run the job in\n#    the foreground.Each job can be identified...\" -> add a space after each\n#    of those two sentence-ending periods.\nReplace-Once \"in the background.Otherwise, it should run the job in the foreground.Each job can\" \"in the background. Otherwise, it should run the job in the foreground. Each job can\"\n\n# 3) \"...manipulating the joblist.)tsh should support the following...\"\n#    -> \"...manipulating the joblist. tsh should support the following...\"\n#    (drops the stray close-paren and adds a space)\nReplace-Once \"manipulating the joblist.)tsh should support\" \"manipulating the joblist. tsh should support\"\n\n# 4) \"...argument can be either a PID or a JID.The fg <job> command...\"\n#    -> \"...argument can be either a PID or a JID. The fg <job> command...\"\nReplace-Once \"argument can be either a PID or a JID.The fg\" \"argument can be either a PID or a JID. The fg\"\n\n# 5) \"...check your work.Reference solution. The Linux executable...\"\n#    -> \"...check your work. Reference solution. The Linux executable...\"\nReplace-Once \"check your work.Reference solution.\" \"check your work. Reference solution.\"\n\n# 6) \"Usage: sdriver.pl [-hv] -t <trace> -s <shellprog> -a <args>&nbsp;\"\n#    -> drop the trailing stray literal \"&nbsp;\" artifact.\nReplace-Once \"-a <args>&nbsp;\" \"-a <args>\"\n\n# 7) The \"Options:\" paragraph starts with a run containing \"  \" followed by a\n#    second run containing the stray literal \"&nbsp;   -v  &nbsp;         \"\n#    prefix before \"Be more verbose ...\". The source edit deletes the first\n#    run entirely and strips the \"&nbsp;   -v  &nbsp;         \" prefix from\n#    the second run, leaving a single run with just \"Be more verbose ...\".\n#    A Find/Replace across the run boundary (removing everything up to and\n#    including the stray prefix) achieves exactly that merge.\nReplace-Once \"  &nbsp;   -v  &nbsp;         Be more verbose\" \"Be more verbose\"\n"}
